# The deck's slide master currently carries the "Integral" design theme
# (dk2/lt2/accent1-6/hlink/folHlink tuned to a green/teal palette). This
# change swaps the live theme's color scheme back to the stock Office
# palette (the "Office Theme" colors), leaving dk1/lt1 (pure black/white)
# untouched. We go through the master's Theme -> ThemeColorScheme so the
# write lands on the theme part that the slide master / every slide layout
# actually renders with.
#
# PowerPoint COM RGB longs are 0x00BBGGRR (little-endian byte order), so
# each target hex color below is converted to that decimal form.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

# index -> (theme slot, target "Office Theme" hex, COM RGB long)
#  1  dk1       000000 -> 0
#  2  lt1       FFFFFF -> 16777215
#  3  dk2       44546A -> 6968388
#  4  lt2       E7E6E6 -> 15132391
#  5  accent1   5B9BD5 -> 13998939
#  6  accent2   ED7D31 -> 3243501
#  7  accent3   A5A5A5 -> 10855845
#  8  accent4   FFC000 -> 49407
#  9  accent5   4472C4 -> 12874308
# 10  accent6   70AD47 -> 4697456
# 11  hlink     0563C1 -> 12673797
# 12  folHlink  954F72 -> 7491477

$colorScheme.Colors(1).RGB = 0
$colorScheme.Colors(2).RGB = 16777215
$colorScheme.Colors(3).RGB = 6968388
$colorScheme.Colors(4).RGB = 15132391
$colorScheme.Colors(5).RGB = 13998939
$colorScheme.Colors(6).RGB = 3243501
$colorScheme.Colors(7).RGB = 10855845
$colorScheme.Colors(8).RGB = 49407
$colorScheme.Colors(9).RGB = 12874308
$colorScheme.Colors(10).RGB = 4697456
$colorScheme.Colors(11).RGB = 12673797
$colorScheme.Colors(12).RGB = 7491477
